# Update gh-pages to output generated at 456a3b4
# Increment the "想去人数" (want-to-go count) column F for a handful of
# rows across the 展览 (sheet1), 本地生活 (sheet3) and 全部类型 (sheet4)
# worksheets.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1792
$ws1.Range("F7").Value = 2249
$ws1.Range("F8").Value = 2169
$ws1.Range("F9").Value = 1139
$ws1.Range("F11").Value = 25
$ws1.Range("F24").Value = 12397
$ws1.Range("F25").Value = 12443
$ws1.Range("F31").Value = 412

# --- 本地生活 sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 109

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1792
$ws4.Range("F8").Value = 2249
$ws4.Range("F9").Value = 2169
$ws4.Range("F10").Value = 1139
$ws4.Range("F12").Value = 109
$ws4.Range("F13").Value = 25
$ws4.Range("F30").Value = 12397
$ws4.Range("F31").Value = 12443
$ws4.Range("F37").Value = 412
